# Weekly Quantity sheet: correct B17, then remove the two rows that no longer
# belong (dates 45095.99999999999 and 45123.99999999999), and correct the
# quantity for the week that used to be row 23 (45130.99999999999).
$wb = $excel.ActiveWorkbook
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")

# Row 17 (2023-06-11 week): 124 -> 88
$wsWeekly.Cells.Item(17, 2).Value = 88

# Delete row 22 (2023-07-16 week, value 184) first so row 18's index is untouched
$wsWeekly.Rows.Item(22).Delete()

# Delete row 18 (2023-06-18 week, value 24)
$wsWeekly.Rows.Item(18).Delete()

# What was row 23 (2023-07-23 week) is now row 21 after the two deletions above;
# its quantity changes from 284 -> 148
$wsWeekly.Cells.Item(21, 2).Value = 148

# Monthly Trend sheet: two quantities revised downward
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Cells.Item(7, 2).Value = 200
$wsMonthly.Cells.Item(8, 2).Value = 400
